# Applies the profit-recalculation refresh captured in the commit's XML diff.
# Each block below corresponds to one changed leve row on one sheet, writing the
# new currentAveragePrice / LevePrice / LeveProfit figures produced by the scheduled
# market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 27405.475
$ws.Range("I21").Value = 19005.666
$ws.Range("K21").Value = 19005.666
$ws.Range("M21").Value = -18537.666

# Row 23
$ws.Range("H23").Value = 27405.475
$ws.Range("I23").Value = 19005.666
$ws.Range("K23").Value = 19005.666
$ws.Range("M23").Value = -18771.666

# Row 34
$ws.Range("H34").Value = 7984.5884
$ws.Range("I34").Value = 5441.385
$ws.Range("J34").Value = 16250
$ws.Range("K34").Value = 5441.385
$ws.Range("L34").Value = 16250
$ws.Range("M34").Value = -5238.385
$ws.Range("N34").Value = -16656

# Row 36
$ws.Range("H36").Value = 7984.5884
$ws.Range("I36").Value = 5441.385
$ws.Range("J36").Value = 16250
$ws.Range("K36").Value = 5441.385
$ws.Range("L36").Value = 16250
$ws.Range("M36").Value = -4726.385
$ws.Range("N36").Value = -17680

# Row 55
$ws.Range("H55").Value = 125329.875
$ws.Range("I55").Value = 166941.5
$ws.Range("K55").Value = 166941.5
$ws.Range("M55").Value = -166727.5

# Row 133
$ws.Range("H133").Value = 16571.428
$ws.Range("J133").Value = 16571.428
$ws.Range("L133").Value = 16571.428
$ws.Range("N133").Value = -26691.428

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1706.375
$ws.Range("I45").Value = 1300.25
$ws.Range("K45").Value = 1300.25
$ws.Range("M45").Value = -923.25

# Row 61
$ws.Range("H61").Value = 2613.862
$ws.Range("I61").Value = 1879.4546
$ws.Range("K61").Value = 1879.4546
$ws.Range("M61").Value = -1667.4546

# Row 74
$ws.Range("H74").Value = 10844.571
$ws.Range("J74").Value = 31700
$ws.Range("L74").Value = 31700
$ws.Range("N74").Value = -33448

# Row 77
$ws.Range("H77").Value = 10844.571
$ws.Range("J77").Value = 31700
$ws.Range("L77").Value = 158500
$ws.Range("N77").Value = -167236

# Row 136
$ws.Range("H136").Value = 2613.862
$ws.Range("I136").Value = 1879.4546
$ws.Range("K136").Value = 5638.3638
$ws.Range("M136").Value = -3088.3638

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 14927.5
$ws.Range("I86").Value = 1901.5
$ws.Range("J86").Value = 27953.5
$ws.Range("K86").Value = 1901.5
$ws.Range("L86").Value = 27953.5
$ws.Range("M86").Value = -778.5
$ws.Range("N86").Value = -30199.5

# Row 89
$ws.Range("H89").Value = 14927.5
$ws.Range("I89").Value = 1901.5
$ws.Range("J89").Value = 27953.5
$ws.Range("K89").Value = 9507.5
$ws.Range("L89").Value = 139767.5
$ws.Range("M89").Value = -3891.5
$ws.Range("N89").Value = -150999.5

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 2675
$ws.Range("J22").Value = 2675
$ws.Range("L22").Value = 8025
$ws.Range("N22").Value = -8363

# Row 27
$ws.Range("H27").Value = 2675
$ws.Range("J27").Value = 2675
$ws.Range("L27").Value = 8025
$ws.Range("N27").Value = -8229

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 24.5
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 28
$ws.Range("K2").Value = 14
$ws.Range("L2").Value = 28
$ws.Range("M2").Value = 99
$ws.Range("N2").Value = -254

# Row 80
$ws.Range("H80").Value = 2294
$ws.Range("I80").Value = 2188.3635
$ws.Range("J80").Value = 2875
$ws.Range("K80").Value = 2188.3635
$ws.Range("L80").Value = 2875
$ws.Range("M80").Value = -1190.3635
$ws.Range("N80").Value = -4871

# Row 83
$ws.Range("H83").Value = 2294
$ws.Range("I83").Value = 2188.3635
$ws.Range("J83").Value = 2875
$ws.Range("K83").Value = 10941.8175
$ws.Range("L83").Value = 14375
$ws.Range("M83").Value = -5949.817499999999
$ws.Range("N83").Value = -24359

$ws = $wb.Worksheets.Item("LTW")
# Row 25
$ws.Range("H25").Value = 35623.75
$ws.Range("I25").Value = 35998
$ws.Range("J25").Value = 35000
$ws.Range("K25").Value = 35998
$ws.Range("L25").Value = 35000
$ws.Range("M25").Value = -35768
$ws.Range("N25").Value = -35460

# Row 26
$ws.Range("H26").Value = 28500

# Row 31
$ws.Range("H31").Value = 4537.375
$ws.Range("I31").Value = 3224.75
$ws.Range("K31").Value = 3224.75
$ws.Range("M31").Value = -2976.75

# Row 46
$ws.Range("H46").Value = 504.73334
$ws.Range("I46").Value = 504.14285
$ws.Range("J46").Value = 505.25
$ws.Range("K46").Value = 504.14285
$ws.Range("L46").Value = 505.25
$ws.Range("M46").Value = -316.14285
$ws.Range("N46").Value = -881.25

# Row 55
$ws.Range("H55").Value = 372.54544
$ws.Range("I55").Value = 385.14285
$ws.Range("J55").Value = 350.5
$ws.Range("K55").Value = 385.14285
$ws.Range("L55").Value = 350.5
$ws.Range("M55").Value = -212.14285
$ws.Range("N55").Value = -696.5

# Row 93
$ws.Range("H93").Value = 628.1852
$ws.Range("I93").Value = 612.0454999999999
$ws.Range("J93").Value = 699.2
$ws.Range("K93").Value = 612.0454999999999
$ws.Range("L93").Value = 699.2
$ws.Range("M93").Value = 635.9545000000001
$ws.Range("N93").Value = -3195.2

# Row 132
$ws.Range("H132").Value = 3672.2856
$ws.Range("I132").Value = 2562.875
$ws.Range("J132").Value = 5151.5
$ws.Range("K132").Value = 7688.625
$ws.Range("L132").Value = 15454.5
$ws.Range("M132").Value = -5158.625
$ws.Range("N132").Value = -20514.5

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 14869.728
$ws.Range("I62").Value = 17709.875
$ws.Range("J62").Value = 7296
$ws.Range("K62").Value = 17709.875
$ws.Range("L62").Value = 7296
$ws.Range("M62").Value = -17085.875
$ws.Range("N62").Value = -8544

# Row 65
$ws.Range("H65").Value = 14869.728
$ws.Range("I65").Value = 17709.875
$ws.Range("J65").Value = 7296
$ws.Range("K65").Value = 88549.375
$ws.Range("L65").Value = 36480
$ws.Range("M65").Value = -85429.375
$ws.Range("N65").Value = -42720

# Row 81
$ws.Range("H81").Value = 3794.44
$ws.Range("I81").Value = 980.1429000000001
$ws.Range("J81").Value = 4888.8887
$ws.Range("K81").Value = 1960.2858
$ws.Range("L81").Value = 9777.777400000001
$ws.Range("M81").Value = -899.2858000000001
$ws.Range("N81").Value = -11899.7774

# Row 84
$ws.Range("H84").Value = 3794.44
$ws.Range("I84").Value = 980.1429000000001
$ws.Range("J84").Value = 4888.8887
$ws.Range("K84").Value = 9801.429
$ws.Range("L84").Value = 48888.887
$ws.Range("M84").Value = -4497.429
$ws.Range("N84").Value = -59496.887

# Row 126
$ws.Range("H126").Value = 46944.953
$ws.Range("I126").Value = 49061.383
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 147184.149
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -144714.149
$ws.Range("N126").Value = -12440

# Row 132
$ws.Range("H132").Value = 11630525
$ws.Range("I132").Value = 14707884
$ws.Range("K132").Value = 44123652
$ws.Range("M132").Value = -44121122
